$d = $word.ActiveDocument

# wdColorRed (255 decimal == 0x0000FF in VBA's BGR wdColor scheme) serialises
# to OOXML as <w:color w:val="FF0000"/>, which is exactly the red used
# throughout this revision.
$RED = 255

# --- Paragraph 1: "Первоначальный вход. Регистрация нового..." ---------
$d.Paragraphs(1).Range.Font.Color = $RED

# --- Paragraph 2: "Убрать. Расскажи о себе." ----------------------------
$d.Paragraphs(2).Range.Font.Color = $RED

# --- Paragraph 3: "Поменять выбор городов." ------------------------------
$d.Paragraphs(3).Range.Font.Color = $RED

# --- Paragraph 4: "Стоп-слова ..." (incl. hyperlink to prnt.sc/8p862...) -
$d.Paragraphs(4).Range.Font.Color = $RED
$d.Hyperlinks(1).Range.Font.Color = $RED

# --- Paragraph 5: "Не прошел регистрацию ..." ----------------------------
# First split the second hyperlink's run into three runs
# ("https://prnt.sc/R_fXO" / "S" / "L9Q-fM") by toggling a trivial
# formatting property on/off around the lone "S" character — this forces a
# run break without altering the final visible formatting, matching the
# run layout seen in the target document. Colour is applied afterwards so
# all three end up identically red.
$hl2 = $d.Hyperlinks(2).Range.Duplicate
$sSplit = $hl2.Duplicate
$sSplit.Find.Execute("S", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sSplit.Bold = 1
$sSplit.Bold = 0

$d.Paragraphs(5).Range.Font.Color = $RED
$d.Hyperlinks(2).Range.Font.Color = $RED

# --- Paragraph 7: "⚙️ Настройка. (то же что и в старом боте)" ------------
$d.Paragraphs(7).Range.Font.Color = $RED

# --- Paragraph 8: "📌 Помощь (то же что и в старом боте)" -----------------
# Merge the " " run with the "(то же что и в старом боте)" run into a
# single run by replacing the combined span with identical text (the
# engine coalesces the matched span into one run when it rewrites it).
$p8 = $d.Paragraphs(8).Range
$p8.Find.Execute(" (то же что и в старом боте)", $true, $false, $false, $false, $false, $true, 1, $false, " (то же что и в старом боте)", 2) | Out-Null
$d.Paragraphs(8).Range.Font.Color = $RED

# --- Paragraph 9: "❤️‍🔥 Совместимость. Вопрос нужна ли в админке ..." ----
# Merge ". " with "Вопрос нужна ли в " into a single run the same way,
# then recolour the whole paragraph (this also turns the old green
# 70AD47/accent6 run colours into the new plain red FF0000).
$p9 = $d.Paragraphs(9).Range
$p9.Find.Execute(". Вопрос нужна ли в ", $true, $false, $false, $false, $false, $true, 1, $false, ". Вопрос нужна ли в ", 2) | Out-Null
$d.Paragraphs(9).Range.Font.Color = $RED
